$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Item #2, "Add APIs"): set Actual Start-date (E4) to 2019-08-22
$ws.Range("E4").NumberFormat = "d-mmm-yy"
$ws.Range("E4").Value = Get-Date -Year 2019 -Month 8 -Day 22 -Hour 0 -Minute 0 -Second 0

# Update Status (G4) from "Not started" to "In-progress" with green font color
$ws.Range("G4").Value = "In-progress"
$ws.Range("G4").Font.Color = 5287936

# Update the active selection to E18
$ws.Range("E18").Select()
